# Updates the cryptos price/volume table (columns D "Price" and E
# "Volume(1h)") with freshly scraped figures, mirroring the daily
# GitHub Actions refresh.
#
# Note: several "Price" values are plain decimal numbers (e.g. "308.77").
# Assigning such a string straight to Range.Value would make Excel parse it
# as a real number, which the source workbook never wants (prices are kept
# as literal text, matching the multi-dot big-number entries like
# "42.282.73" that cannot be numbers at all). So for those cells we assign
# with a leading apostrophe - the normal Excel way to force text entry -
# and then reset .Style back to "Normal" so the cell doesn't keep the
# quote-prefix formatting flag.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.282.73'
$ws.Range('D3').Value = '2.270.10'
$ws.Range('E3').Value = '  -0.45%  '
$ws.Range('E4').Value = '  +0.00%  '
$c = $ws.Range('D5')
$c.Value = "'308.77"
$c.Style = 'Normal'
$ws.Range('E5').Value = '  +0.60%  '
$c = $ws.Range('D6')
$c.Value = "'97.46"
$c.Style = 'Normal'
$ws.Range('E6').Value = '  -0.47%  '
$ws.Range('E7').Value = '  -0.81%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('E9').Value = '  -1.18%  '
$c = $ws.Range('D10')
$c.Value = "'34.99"
$c.Style = 'Normal'
$ws.Range('E10').Value = '  -2.72%  '
$c = $ws.Range('D11')
$c.Value = "'0.0811"
$c.Style = 'Normal'
$ws.Range('E11').Value = '  +1.25%  '
$ws.Range('E12').Value = '  +1.12%  '
$c = $ws.Range('D13')
$c.Value = "'6.85"
$c.Style = 'Normal'
$ws.Range('E13').Value = '  +2.15%  '
$ws.Range('D14').Value = '2.622.54'
$ws.Range('E14').Value = '  -0.42%  '
$c = $ws.Range('D15')
$c.Value = "'14.57"
$c.Style = 'Normal'
$ws.Range('E15').Value = '  +0.73%  '
$ws.Range('D16').Value = '2.263.03'
$ws.Range('E16').Value = '  -1.01%  '
$ws.Range('E17').Value = '  -1.46%  '
$ws.Range('D18').Value = '42.178.16'
$ws.Range('E18').Value = '  -0.26%  '
$ws.Range('E19').Value = '  -2.59%  '
$ws.Range('D20').Value = '0.0₃0905'
$ws.Range('E20').Value = '  -0.90%  '
$ws.Range('E21').Value = '  -1.02%  '
$c = $ws.Range('D22')
$c.Value = "'67.63"
$c.Style = 'Normal'
$ws.Range('E22').Value = '  -0.12%  '
$c = $ws.Range('D23')
$c.Value = "'236.68"
$c.Style = 'Normal'
$ws.Range('E23').Value = '  -2.18%  '
$c = $ws.Range('D24')
$c.Value = "'2.59"
$c.Style = 'Normal'
$ws.Range('E24').Value = '  -0.61%  '
$c = $ws.Range('D25')
$c.Value = "'1.98"
$c.Style = 'Normal'
$ws.Range('E25').Value = '  +1.63%  '
$ws.Range('E26').Value = '  +0.04%  '
$c = $ws.Range('D27')
$c.Value = "'23.61"
$c.Style = 'Normal'
$ws.Range('E27').Value = '  -1.22%  '
$c = $ws.Range('D28')
$c.Value = "'37.09"
$c.Style = 'Normal'
$ws.Range('E28').Value = '  -2.20%  '
$ws.Range('E29').Value = '  +0.05%  '
$ws.Range('E30').Value = '  +0.61%  '
$c = $ws.Range('D31')
$c.Value = "'163.47"
$c.Style = 'Normal'
$ws.Range('E31').Value = '  +1.71%  '
$ws.Range('E32').Value = '  -0.19%  '
$ws.Range('E33').Value = '  +0.06%  '
$ws.Range('E34').Value = '  -0.86%  '
$c = $ws.Range('D35')
$c.Value = "'17.57"
$c.Style = 'Normal'
$ws.Range('E35').Value = '  +2.28%  '
$ws.Range('E36').Value = '  -1.80%  '
$ws.Range('E37').Value = '  -0.08%  '
$ws.Range('E38').Value = '  -2.69%  '
$ws.Range('E39').Value = '  -0.32%  '
$ws.Range('E40').Value = '  -2.06%  '
$c = $ws.Range('D41')
$c.Value = "'4.15"
$c.Style = 'Normal'
$ws.Range('E41').Value = '  +0.01%  '
$ws.Range('E42').Value = '  -5.04%  '
$ws.Range('D43').Value = '1.948.10'
$ws.Range('E43').Value = '  -2.75%  '
$ws.Range('E44').Value = '  -1.26%  '
$c = $ws.Range('D45')
$c.Value = "'18.71"
$c.Style = 'Normal'
$ws.Range('E45').Value = '  -1.83%  '
$ws.Range('E46').Value = '  -2.21%  '
$c = $ws.Range('D47')
$c.Value = "'9.78"
$c.Style = 'Normal'
$ws.Range('E47').Value = '  -2.74%  '
$c = $ws.Range('D48')
$c.Value = "'54.46"
$c.Style = 'Normal'
$ws.Range('E48').Value = '  +1.91%  '
$ws.Range('D49').Value = '2.494.61'
$ws.Range('E49').Value = '  -0.43%  '
$ws.Range('E50').Value = '  -0.86%  '
$c = $ws.Range('D51')
$c.Value = "'71.51"
$c.Style = 'Normal'
$ws.Range('E51').Value = '  -1.26%  '
